# Updated cryptos list values (price + volume%) per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.931.51"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "3.477.50"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'601.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "'142.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").Value = "3.476.09"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").Value = "'8.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.85%  "
$ws.Range("E11").Value = "  -4.59%  "
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "4.073.55"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'0.0000202"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "'30.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "3.478.30"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "66.120.46"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'10.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("D20").Value = "'6.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").Value = "'14.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "'419.29"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").Value = "'0.588"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.39%  "
$ws.Range("D24").Value = "'77.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "3.620.45"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'0.0000114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.23%  "
$ws.Range("D28").Value = "'9.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.72%  "
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("E30").Value = "  -1.07%  "
$ws.Range("D31").Value = "'1.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("E33").Value = "  -6.89%  "
$ws.Range("D34").Value = "'25.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "3.480.26"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("D38").Value = "'5.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.69%  "
$ws.Range("D39").Value = "'7.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "'169.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "'0.0864"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.88%  "
$ws.Range("D43").Value = "'0.888"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "'5.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.32%  "
$ws.Range("E45").Value = "  -7.01%  "
$ws.Range("D46").Value = "'45.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("D47").Value = "'25.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.55%  "
$ws.Range("D48").Value = "'1.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.20%  "
$ws.Range("D51").Value = "'0.927"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.68%  "
